# Updates cryptos list values (Price / Volume(1h) columns, plus one rank-51 coin swap)
# to match the latest scraped data, mirroring the upstream GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''63.689.18'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.92%  '
# Row 3
$ws.Range('D3').Value = '''3.287.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.74%  '
# Row 4
$ws.Range('E4').Value = '  +0.10%  '
# Row 5
$ws.Range('D5').Value = '''604.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.67%  '
# Row 6
$ws.Range('D6').Value = '''141.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.10%  '
# Row 7
$ws.Range('E7').Value = '  +0.02%  '
# Row 8
$ws.Range('D8').Value = '''3.285.57'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.74%  '
# Row 9
$ws.Range('D9').Value = '''0.518'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.42%  '
# Row 10
$ws.Range('E10').Value = '  +2.61%  '
# Row 11
$ws.Range('D11').Value = '''5.44'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.47%  '
# Row 12
$ws.Range('E12').Value = '  +2.52%  '
# Row 13
$ws.Range('D13').Value = '''0.0000246'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.50%  '
# Row 14
$ws.Range('D14').Value = '''34.43'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.92%  '
# Row 15
$ws.Range('D15').Value = '''3.838.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.03%  '
# Row 16
$ws.Range('E16').Value = '  +0.98%  '
# Row 17
$ws.Range('D17').Value = '''3.291.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.91%  '
# Row 18
$ws.Range('D18').Value = '''63.762.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.03%  '
# Row 19
$ws.Range('D19').Value = '''6.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.35%  '
# Row 20
$ws.Range('D20').Value = '''479.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.58%  '
# Row 21
$ws.Range('D21').Value = '''14.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.60%  '
# Row 22
$ws.Range('E22').Value = '  +4.15%  '
# Row 23
$ws.Range('D23').Value = '''8.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.50%  '
# Row 24
$ws.Range('D24').Value = '''13.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.52%  '
# Row 25
$ws.Range('D25').Value = '''84.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.49%  '
# Row 27
$ws.Range('E27').Value = '  +1.81%  '
# Row 28
$ws.Range('E28').Value = '  +0.08%  '
# Row 29
$ws.Range('E29').Value = '  +5.60%  '
# Row 30
$ws.Range('E30').Value = '  +2.09%  '
# Row 31
$ws.Range('E31').Value = '  +3.09%  '
# Row 32
$ws.Range('D32').Value = '''28.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.37%  '
# Row 33
$ws.Range('E33').Value = '  -2.45%  '
# Row 34
$ws.Range('E34').Value = '  -0.21%  '
# Row 35
$ws.Range('E35').Value = '  +2.77%  '
# Row 36
$ws.Range('D36').Value = '''5.95'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.34%  '
# Row 37
$ws.Range('D37').Value = '''53.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.57%  '
# Row 38
$ws.Range('D38').Value = '''0.0₃0736'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.75%  '
# Row 39
$ws.Range('E39').Value = '  +2.71%  '
# Row 40
$ws.Range('D40').Value = '''427.88'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.82%  '
# Row 41
$ws.Range('D41').Value = '''3.062.92'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.36%  '
# Row 42
$ws.Range('E42').Value = '  +1.44%  '
# Row 43
$ws.Range('D43').Value = '''2.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.92%  '
# Row 44
$ws.Range('E44').Value = '  +0.25%  '
# Row 45
$ws.Range('E45').Value = '  +1.02%  '
# Row 46
$ws.Range('E46').Value = '  +2.94%  '
# Row 47
$ws.Range('E47').Value = '  +0.04%  '
# Row 48
$ws.Range('D48').Value = '''26.09'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.68%  '
# Row 49
$ws.Range('E49').Value = '  +1.39%  '
# Row 50
$ws.Range('D50').Value = '''125.41'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.88%  '
# Row 51
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').Value = '''35.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +12.50%  '
